# Update current CRISPR template with library names and pipeline names

$wb = $excel.ActiveWorkbook

# --- GUIDES sheet: add "SUMO" and "NEW_LIBRARY" entries, drop the retired
#     pLib017-reporter / pLib017-insert entries, and move the two compound
#     entries (pLib013-pLib014, NGN-NGG) to the bottom of the list.
$guides = $wb.Worksheets.Item("GUIDES")
$guides.Unprotect()
$guides.Range("A1").Value = "TKOv3"
$guides.Range("A2").Value = "NGN"
$guides.Range("A3").Value = "NGG"
$guides.Range("A4").Value = "pLib013"
$guides.Range("A5").Value = "pLib014"
$guides.Range("A6").Value = "pLib016"
$guides.Range("A7").Value = "TLS_ABE"
$guides.Range("A8").Value = "TOP2_ABE"
$guides.Range("A9").Value = "pLib017"
$guides.Range("A10").Value = "SUMO"
$guides.Range("A11").Value = "NEW_LIBRARY"
$guides.Range("A12").Value = "pLib013-pLib014"
$guides.Range("A13").Value = "NGN-NGG"
$guides.Range("A4").Select()

# --- pipeline sheet: add two new pipeline options between MAGECK and
#     CRISPRESSO2, and widen column A to fit the longer names.
$pipeline = $wb.Worksheets.Item("pipeline")
$pipeline.Unprotect()
$pipeline.Range("A1").Value = "MAGECK"
$pipeline.Range("A2").Value = "MAGECK-DRUGZ-BAGEL"
$pipeline.Range("A3").Value = "MAGECK-BEAN"
$pipeline.Range("A4").Value = "CRISPRESSO2"
$pipeline.Range("A1").ColumnWidth = 20.6
$pipeline.Range("A4").Select()

# --- samples sheet: leave data as-is, just move the selection / keep it
#     the active sheet/tab.
$samples = $wb.Worksheets.Item("samples")
$samples.Activate()
$samples.Range("G1").Select()
